$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.319.48"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.874.66"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "0.7102"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").Value = "242.01"

$ws.Range("D8").Value = "0.07803"
$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("D9").Value = "0.3107"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "25.14"
$ws.Range("E10").Value = "  +1.57%  "

$ws.Range("D11").Value = "0.08424"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "1.862.15"
$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").Value = "5.237"
$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "0.7157"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("D15").Value = "91.17"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Value = "29.327.03"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").Value = "0.000008324"
$ws.Range("E17").Value = "  +1.65%  "

$ws.Range("D18").Value = "6.083"
$ws.Range("E18").Value = "  +2.43%  "

$ws.Range("D19").Value = "240.54"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Value = "2.115.85"
$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "7.746"

$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "0.1597"
$ws.Range("E25").Value = "  -2.08%  "

$ws.Range("D26").Value = "162.40"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").Value = "9.032"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "18.50"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "1.506"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("D30").Value = "4.404"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").Value = "4.319"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("D33").Value = "0.05378"
$ws.Range("E33").Value = "  +3.95%  "

$ws.Range("E34").Value = "  +1.75%  "

$ws.Range("E35").Value = "  +0.64%  "

$ws.Range("D36").Value = "0.7495"
$ws.Range("E36").Value = "  -3.10%  "

$ws.Range("D37").Value = "2.690"
$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "0.01884"
$ws.Range("E38").Value = "  +1.28%  "

$ws.Range("D39").Value = "1.232.92"
$ws.Range("E39").Value = "  +6.31%  "

$ws.Range("D40").Value = "2.728"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").Value = "6.482"
$ws.Range("E41").Value = "  +1.15%  "

$ws.Range("D42").Value = "0.8954"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("D43").Value = "109.18"
$ws.Range("E43").Value = "  +4.25%  "

$ws.Range("D44").Value = "72.39"
$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "2.020.35"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.797"
$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.5200"
$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.451"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.4335"
$ws.Range("E51").Value = "  +0.84%  "
